# "ein Teil von Reliable geht" - add the statistics header row (Testart,
# Datum, Zeit, Anzahl Client-Threads, Durchschnittliche RTT, Maximale RTT,
# Minimale RTT) to the first worksheet (Tabelle1 / Stats.xlsx), size the
# first three columns, and leave the selection on G1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A1").Value = "Testart"
$ws.Range("B1").Value = "Datum"
$ws.Range("C1").Value = "Zeit"
$ws.Range("D1").Value = "Anzahl Client-Threads"
$ws.Range("E1").Value = "Durchschnittliche RTT"
$ws.Range("F1").Value = "Maximale RTT"
$ws.Range("G1").Value = "Minimale RTT "

# Column widths (characters) matching the saved workbook's <cols> entries.
$ws.Columns.Item(1).ColumnWidth = 11.166666666666666
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666

# Leave the active selection on G1, as in the committed workbook.
$ws.Range("G1").Select()
